$wb = $excel.ActiveWorkbook

$oneWay = $wb.Worksheets.Item("Air_Sabre_OneWay")
$roundTrip = $wb.Worksheets.Item("Air_Sabre_RoundTrip")

$newPipeline = "LOGIN|Search|AddToCart|CHECKOUTTRIP|ENTERPAXINFO|CONFIRMPAXINFO"

# Air_Sabre_OneWay sheet
$oneWay.Range("B2").Value = $newPipeline
$oneWay.Range("B3").Value = $newPipeline
$oneWay.Range("B4").Value = $newPipeline
$oneWay.Range("B5").Value = $newPipeline

$oneWay.Range("D2").Value = "ATL-FLL"
$oneWay.Range("D3").Value = "LAX-LAS"
$oneWay.Range("D4").Value = "ATL-FLL"
$oneWay.Range("D5").Value = "ATL-FLL"

# Air_Sabre_RoundTrip sheet
$roundTrip.Range("B2").Value = $newPipeline
$roundTrip.Range("B3").Value = $newPipeline
$roundTrip.Range("B4").Value = $newPipeline
$roundTrip.Range("B5").Value = $newPipeline

$roundTrip.Range("D2").Value = "ATL-FLL"
$roundTrip.Range("D3").Value = "LAX-LAS"
$roundTrip.Range("D4").Value = "ATL-FLL"
$roundTrip.Range("D5").Value = "ATL-FLL"

# Selections / view state as per diff
$oneWay.Range("D2:D5").Select()
$oneWay.Application.ActiveWindow.ScrollRow = 3

$roundTrip.Range("B2").Select()
